# Adds the missing "Média" (average) column values to the
# "Tempo Total(Segundos)" table, and corrects the previously-wrong
# "Média" column values (a stale copy-paste of 63,17) in the
# "Taxa de compressão" table.

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Table 3: "Tempo Total(Segundos)" - the last column ("Média") is
# currently empty for the egg/landscape/pattern/zebra rows; fill it
# in with the row averages and bottom-align the cell, matching the
# formatting already used by the other data cells in these rows.
# ---------------------------------------------------------------

function Set-EmptyMediaCell {
    param(
        [int]$TableIndex,
        [int]$RowIndex,
        [int]$ColIndex,
        [string]$Value
    )

    $t = $d.Tables.Item($TableIndex)
    $cell = $t.Cell($RowIndex, $ColIndex)
    $cell.VerticalAlignment = 3   # wdCellAlignVerticalBottom
    $rng = $cell.Range
    $rng.Text = $Value

    $t = $d.Tables.Item($TableIndex)
    $cell = $t.Cell($RowIndex, $ColIndex)
    $rng = $cell.Range
    $rng.Font.Size = 8

    $t = $d.Tables.Item($TableIndex)
    $cell = $t.Cell($RowIndex, $ColIndex)
    $rng = $cell.Range
    $rng.Font.SizeBi = 8

    $t = $d.Tables.Item($TableIndex)
    $cell = $t.Cell($RowIndex, $ColIndex)
    $rng = $cell.Range
    $rng.Font.Color = 0

    $t = $d.Tables.Item($TableIndex)
    $cell = $t.Cell($RowIndex, $ColIndex)
    $rng = $cell.Range
    $find = $rng.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Replacement.Font.Name = "Calibri"
    $find.Replacement.Font.NameBi = "Calibri"
    $find.Execute($Value, $false, $false, $false, $false, $false, $true, 1, $false, $Value, 2) | Out-Null
}

Set-EmptyMediaCell 3 3 17 "73,24"   # egg.bmp
Set-EmptyMediaCell 3 4 17 "50,74"   # landscape.bmp
Set-EmptyMediaCell 3 5 17 "148,14"  # pattern.bmp
Set-EmptyMediaCell 3 6 17 "74,99"   # zebra.bmp

# ---------------------------------------------------------------
# Table 5: "Taxa de compressão" - the "Média" column had been filled
# with a stale "63,17" placeholder for three of the four rows; swap
# in the real averages (the zebra.bmp row's average already was
# 63,17, so it is left untouched).
# ---------------------------------------------------------------

$t5 = $d.Tables.Item(5)
$t5.Cell(3, 17).Range.Text = "69,65"   # egg.bmp

$t5 = $d.Tables.Item(5)
$t5.Cell(4, 17).Range.Text = "67,74"   # landscape.bmp

$t5 = $d.Tables.Item(5)
$t5.Cell(5, 17).Range.Text = "94,4"    # pattern.bmp
